$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 197 - this shifts the existing rows
# 197..284 down to 198..285 and extends the used range accordingly.
$ws.Rows("197:197").Insert()

# Populate the newly inserted row 197 with the new weekly data point.
$ws.Range("A197").Value = 8
$ws.Range("B197").Value = "Terminal La Palmera de La Serena"
$ws.Range("C197").Value = "Coquimbo"
$ws.Range("D197").Value = 44960
$ws.Range("E197").Value = 4
$ws.Range("F197").Value = 100112037
$ws.Range("G197").Value = "Cebollín"
$ws.Range("H197").Value = "Sin especificar"
$ws.Range("I197").Value = "Primera"
$ws.Range("J197").Value = 1300
$ws.Range("K197").Value = 1200
$ws.Range("L197").Value = 1400
$ws.Range("M197").Value = 1300
$ws.Range("N197").Value = "$/paquete 6 unidades"
$ws.Range("O197").Value = "Provincia del Elquí"
$ws.Range("P197").Value = 217
$ws.Range("Q197").Value = 6
$ws.Range("R197").Value = "Hortaliza"
